$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - plain (non-shared) formulas, A4 computed via formula 1/5
$ws.Range("A4").Formula = "=1/5"
$ws.Range("B4").Value = "AMX"
$ws.Range("C4").Value = 58
$ws.Range("D4").Value = 26
$ws.Range("E4").Value = 22
$ws.Range("F4").Value = 113
$ws.Range("G4").Formula = "=F4/(F4+E4)"
$ws.Range("H4").Formula = "=C4/(C4+D4)"
$ws.Range("I4").Formula = "=1-H4"

# Row 5 - plain (non-shared) formulas, A5 is a literal value
$ws.Range("A5").Value = 0.1667
$ws.Range("B5").Value = "AMX"
$ws.Range("C5").Value = 58
$ws.Range("D5").Value = 26
$ws.Range("E5").Value = 22
$ws.Range("F5").Value = 113
$ws.Range("G5").Formula = "=F5/(F5+E5)"
$ws.Range("H5").Formula = "=C5/(C5+D5)"
$ws.Range("I5").Formula = "=1-H5"

# Row 6 - literal value, master of shared formula group (rows 6:7)
$ws.Range("A6").Value = 0.1429
$ws.Range("B6").Value = "AMX"
$ws.Range("C6").Value = 58
$ws.Range("D6").Value = 26
$ws.Range("E6").Value = 22
$ws.Range("F6").Value = 113

# Row 7 - literal value
$ws.Range("A7").Value = 0.5
$ws.Range("B7").Value = "AMX"
$ws.Range("C7").Value = 58
$ws.Range("D7").Value = 26
$ws.Range("E7").Value = 22
$ws.Range("F7").Value = 113

# Fill the shared formula groups across rows 6:7 in one shot so that
# Excel stores them as shared formulas (t="shared")
$ws.Range("G6:G7").Formula = "=F6/(F6+E6)"
$ws.Range("H6:H7").Formula = "=C6/(C6+D6)"
$ws.Range("I6:I7").Formula = "=1-H6"

# Update the selection to match the post-edit state
$ws.Range("I5:I7").Select()
